# Generate Report for Handoff
# Adds a new file entry (dd233af3-56c7-4c64-a07e-cf855ca0f82d) as row 8 to the
# Overview, zh-cn and de-de worksheets, mirroring the existing rows' layout.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: Overview
# ---------------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

$ov.Range("B8").Value = "Ready for handoff"
$ov.Range("C8").Value = "Ready for handoff"
$ov.Range("D8").Value = "2016-45-18 02:45:02"

$ov.Hyperlinks.Add(
    $ov.Range("A8"),
    "https://github.com/OpenLocalizationTest/oltest/blob/dd233af356c74c64a07ecf855ca0f82d0000000/e2e/dd233af3-56c7-4c64-a07e-cf855ca0f82d.md",
    "",
    "",
    "dd233af3-56c7-4c64-a07e-cf855ca0f82d.md"
) | Out-Null

# ---------------------------------------------------------------------------
# Sheet 2: zh-cn
# ---------------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("C8").Value = "Ready for handoff"
$zh.Range("E8").Value = "2016-03-18 02:44:53"
$zh.Range("H8").Value = "0001-01-01 00:00:00"
$zh.Range("I8").Value = "Include"

$zh.Hyperlinks.Add(
    $zh.Range("A8"),
    "https://github.com/OpenLocalizationTest/oltest/blob/dd233af356c74c64a07ecf855ca0f82d0000000/e2e/dd233af3-56c7-4c64-a07e-cf855ca0f82d.md",
    "",
    "",
    "dd233af3-56c7-4c64-a07e-cf855ca0f82d.md"
) | Out-Null

$zh.Hyperlinks.Add(
    $zh.Range("B8"),
    "https://github.com/OpenLocalizationTest/oltest/blob/dd233af356c74c64a07ecf855ca0f82d0000000/e2e/dd233af3-56c7-4c64-a07e-cf855ca0f82d.md",
    "",
    "",
    ".md"
) | Out-Null

$zh.Hyperlinks.Add(
    $zh.Range("D8"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0000000000000000000000000000000000000000/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/xinjiang/ht/dd233af3-56c7-4c64-a07e-cf855ca0f82d.a3e3d5d1a309f9944f6103f007df5fd5c012a303.zh-cn.xlf",
    "",
    "",
    "dd233af3-56c7-4c64-a07e-cf855ca0f82d.a3e3d5d1a309f9944f6103f007df5fd5c012a303.zh-cn.xlf"
) | Out-Null

# ---------------------------------------------------------------------------
# Sheet 3: de-de
# ---------------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Range("C8").Value = "Ready for handoff"
$de.Range("E8").Value = "2016-03-18 02:45:02"
$de.Range("H8").Value = "0001-01-01 00:00:00"
$de.Range("I8").Value = "Include"

$de.Hyperlinks.Add(
    $de.Range("A8"),
    "https://github.com/OpenLocalizationTest/oltest/blob/dd233af356c74c64a07ecf855ca0f82d0000000/e2e/dd233af3-56c7-4c64-a07e-cf855ca0f82d.md",
    "",
    "",
    "dd233af3-56c7-4c64-a07e-cf855ca0f82d.md"
) | Out-Null

$de.Hyperlinks.Add(
    $de.Range("B8"),
    "https://github.com/OpenLocalizationTest/oltest/blob/dd233af356c74c64a07ecf855ca0f82d0000000/e2e/dd233af3-56c7-4c64-a07e-cf855ca0f82d.md",
    "",
    "",
    ".md"
) | Out-Null

$de.Hyperlinks.Add(
    $de.Range("D8"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0000000000000000000000000000000000000000/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/xinjiang/ht/dd233af3-56c7-4c64-a07e-cf855ca0f82d.a3e3d5d1a309f9944f6103f007df5fd5c012a303.de-de.xlf",
    "",
    "",
    "dd233af3-56c7-4c64-a07e-cf855ca0f82d.a3e3d5d1a309f9944f6103f007df5fd5c012a303.de-de.xlf"
) | Out-Null
